$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "USERNAME"
$ws.Range("B8").Value = "PASSWORD"
$ws.Range("C8").Value = 1

$ws.Range("A8:C8").HorizontalAlignment = -4108

$ws.Range("E20").Select()
